$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl11"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8417533333333332
$ws.Range("H2").Value = 2.52526
$ws.Range("I2").Value = 0.01079423211523897
$ws.Range("J2").Value = 0.01079423211523897
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.227228
$ws.Range("N2").Value = 0.681684
$ws.Range("O2").Value = 0.2376267857721762
$ws.Range("P2").Value = 0.2376267857721762
$ws.Range("Q2").Value = 0.1912699264266666
$ws.Range("R2").Value = 1.72142933784
$ws.Range("S2").Value = 0.002564998682423035
$ws.Range("T2").Value = 0.002564998682423034

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl11"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8417533333333332
$ws.Range("H3").Value = 2.52526
$ws.Range("I3").Value = 0.01079423211523897
$ws.Range("J3").Value = 0.01079423211523897
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6537306666666667
$ws.Range("N3").Value = 1.961192
$ws.Range("O3").Value = 0.6836477770376096
$ws.Range("P3").Value = 0.6836477770376095
$ws.Range("Q3").Value = 0.5502799677688889
$ws.Range("R3").Value = 4.95251970992
$ws.Range("S3").Value = 0.007379452790411095
$ws.Range("T3").Value = 0.007379452790411091

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl11"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8417533333333332
$ws.Range("H4").Value = 2.52526
$ws.Range("I4").Value = 0.01079423211523897
$ws.Range("J4").Value = 0.01079423211523897
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07528033333333332
$ws.Range("N4").Value = 0.225841
$ws.Range("O4").Value = 0.0787254371902143
$ws.Range("P4").Value = 0.0787254371902143
$ws.Range("Q4").Value = 0.06336747151777776
$ws.Range("R4").Value = 0.57030724366
$ws.Range("S4").Value = 0.0008497806424048395
$ws.Range("T4").Value = 0.0008497806424048391

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl11"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 69.05064766666666
$ws.Range("H5").Value = 207.151943
$ws.Range("I5").Value = 0.885471656726338
$ws.Range("J5").Value = 0.8854716567263378
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.227228
$ws.Range("N5").Value = 0.681684
$ws.Range("O5").Value = 0.2376267857721762
$ws.Range("P5").Value = 0.2376267857721762
$ws.Range("Q5").Value = 15.69024056800133
$ws.Range("R5").Value = 141.212165112012
$ws.Range("S5").Value = 0.2104117836802435
$ws.Range("T5").Value = 0.2104117836802434

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl11"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 69.05064766666666
$ws.Range("H6").Value = 207.151943
$ws.Range("I6").Value = 0.885471656726338
$ws.Range("J6").Value = 0.8854716567263378
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6537306666666667
$ws.Range("N6").Value = 1.961192
$ws.Range("O6").Value = 0.6836477770376096
$ws.Range("P6").Value = 0.6836477770376095
$ws.Range("Q6").Value = 45.14052593289511
$ws.Range("R6").Value = 406.264733396056
$ws.Range("S6").Value = 0.6053507297507703
$ws.Range("T6").Value = 0.6053507297507701

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl11"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 69.05064766666666
$ws.Range("H7").Value = 207.151943
$ws.Range("I7").Value = 0.885471656726338
$ws.Range("J7").Value = 0.8854716567263378
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07528033333333332
$ws.Range("N7").Value = 0.225841
$ws.Range("O7").Value = 0.0787254371902143
$ws.Range("P7").Value = 0.0787254371902143
$ws.Range("Q7").Value = 5.198155773229221
$ws.Range("R7").Value = 46.78340195906299
$ws.Range("S7").Value = 0.06970914329532432
$ws.Range("T7").Value = 0.06970914329532429

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ccl11"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.697976666666666
$ws.Range("H8").Value = 23.09393
$ws.Range("I8").Value = 0.09871507918910555
$ws.Range("J8").Value = 0.09871507918910553
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.227228
$ws.Range("N8").Value = 0.681684
$ws.Range("O8").Value = 0.2376267857721762
$ws.Range("P8").Value = 0.2376267857721762
$ws.Range("Q8").Value = 1.749195842013333
$ws.Range("R8").Value = 15.74276257812
$ws.Range("S8").Value = 0.02345734697495299
$ws.Range("T8").Value = 0.02345734697495299

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ccl11"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.697976666666666
$ws.Range("H9").Value = 23.09393
$ws.Range("I9").Value = 0.09871507918910555
$ws.Range("J9").Value = 0.09871507918910553
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6537306666666667
$ws.Range("N9").Value = 1.961192
$ws.Range("O9").Value = 0.6836477770376096
$ws.Range("P9").Value = 0.6836477770376095
$ws.Range("Q9").Value = 5.032403418284445
$ws.Range("R9").Value = 45.29163076456
$ws.Range("S9").Value = 0.0674863444477236
$ws.Range("T9").Value = 0.06748634444772358

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ccl11"
$ws.Range("C10").Value = "Ackr4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.697976666666666
$ws.Range("H10").Value = 23.09393
$ws.Range("I10").Value = 0.09871507918910555
$ws.Range("J10").Value = 0.09871507918910553
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.07528033333333332
$ws.Range("N10").Value = 0.225841
$ws.Range("O10").Value = 0.0787254371902143
$ws.Range("P10").Value = 0.0787254371902143
$ws.Range("Q10").Value = 0.5795062494588888
$ws.Range("R10").Value = 5.21555624513
$ws.Range("S10").Value = 0.00777138776642896
$ws.Range("T10").Value = 0.007771387766428958

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Ccl11"
$ws.Range("C11").Value = "Ackr4"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.391393
$ws.Range("H11").Value = 1.174179
$ws.Range("I11").Value = 0.005019031969317685
$ws.Range("J11").Value = 0.005019031969317684
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.227228
$ws.Range("N11").Value = 0.681684
$ws.Range("O11").Value = 0.2376267857721762
$ws.Range("P11").Value = 0.2376267857721762
$ws.Range("Q11").Value = 0.088935448604
$ws.Range("R11").Value = 0.800419037436
$ws.Range("S11").Value = 0.001192656434556757
$ws.Range("T11").Value = 0.001192656434556757

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Ccl11"
$ws.Range("C12").Value = "Ackr4"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.391393
$ws.Range("H12").Value = 1.174179
$ws.Range("I12").Value = 0.005019031969317685
$ws.Range("J12").Value = 0.005019031969317684
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6537306666666667
$ws.Range("N12").Value = 1.961192
$ws.Range("O12").Value = 0.6836477770376096
$ws.Range("P12").Value = 0.6836477770376095
$ws.Range("Q12").Value = 0.2558656068186667
$ws.Range("R12").Value = 2.302790461368
$ws.Range("S12").Value = 0.003431250048704731
$ws.Range("T12").Value = 0.00343125004870473

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Ccl11"
$ws.Range("C13").Value = "Ackr4"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.391393
$ws.Range("H13").Value = 1.174179
$ws.Range("I13").Value = 0.005019031969317685
$ws.Range("J13").Value = 0.005019031969317684
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.07528033333333332
$ws.Range("N13").Value = 0.225841
$ws.Range("O13").Value = 0.0787254371902143
$ws.Range("P13").Value = 0.0787254371902143
$ws.Range("Q13").Value = 0.02946419550433333
$ws.Range("R13").Value = 0.265177759539
$ws.Range("S13").Value = 0.000395125486056197
$ws.Range("T13").Value = 0.0003951254860561969

$ws.Rows("14:17").Delete()
